# Remove stray "Sheet" row (row 16) from the optimization_parameters sheet,
# then make threshold_b the active sheet (matching the author's final
# view state in the workbook).

$wb = $excel.ActiveWorkbook

$wsOpt = $wb.Worksheets.Item("optimization_parameters")
$wsOpt.Rows.Item(16).Delete()

$wsThreshold = $wb.Worksheets.Item("threshold_b")
$wsThreshold.Activate()
$wsThreshold.Range("A2").Select()
